# Add a new "2022-Q4" worksheet (with its own fund table) right before the
# existing "2022-Q2" sheet, and insert a matching summary row into the
# "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q4" sheet by duplicating the "2022-Q2" sheet so
#    that it inherits the same column widths / cell styles, then replace
#    its contents with the 2022-Q4 fund data.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Drop the extra rows that came from the 2022-Q2 sheet (it had 10 fund
# rows, 2022-Q4 only has 2).
$q4.Range("A4:H11").Delete(-4162)

# Row 1 (headers) is already correct, copied verbatim from 2022-Q2
# ("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比",
#  "持有市值(亿元)", "仓位排名").

# Force the numeric-looking fund figures (B,D,E,F,G) to stay plain text
# (matching the source data, which stores them as text, not numbers) by
# applying a text number format before assigning, then clearing the
# format afterwards so no stray style is left behind.
$textRange = $q4.Range("B2:G3")
$textRange.NumberFormat = "@"

$q4.Range("B2").Value = "008861"
$q4.Range("C2").Value = "西部利得港股通新机遇灵活配置混合A"
$q4.Range("D2").Value = "0.25"
$q4.Range("E2").Value = "87.69"
$q4.Range("F2").Value = "3.46"
$q4.Range("G2").Value = "0.0086"
$q4.Range("H2").Value = 8

$q4.Range("B3").Value = "010093"
$q4.Range("C3").Value = "西部利得港股通新机遇灵活配置混合C"
$q4.Range("D3").Value = "0.12"
$q4.Range("E3").Value = "87.69"
$q4.Range("F3").Value = "3.46"
$q4.Range("G3").Value = "0.0042"
$q4.Range("H3").Value = 8

$textRange.ClearFormats()

# Index column (A2/A3) keeps the same style it already had from the
# copied sheet (s="2"), just refresh the values.
$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1

# ---------------------------------------------------------------------
# 2) Insert a new row into the "总计" sheet for 2022-Q4 and reindex the
#    summary rows (0, 1, 2).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Restore formatting for the inserted row's index cell (A2) so it keeps
# the same style as the other index cells, then clear the stray style
# that Insert() copied onto B2:D2.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
